# Auto-generated script to apply 2022-05-15 violent crime data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 2266
$ws.Range('I3').Value = 2427
$ws.Range('G4').Value = 1429
$ws.Range('H4').Value = 1661
$ws.Range('I4').Value = 592
$ws.Range('I6').Value = 2818
$ws.Range('G7').Value = 24652
$ws.Range('H7').Value = 25971
$ws.Range('I7').Value = 8317

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I3').Value = 18
$ws.Range('I7').Value = 87

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('I6').Value = 20
$ws.Range('I7').Value = 41

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('I6').Value = 7
$ws.Range('I7').Value = 26

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I2').Value = 83
$ws.Range('I7').Value = 265

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I2').Value = 39
$ws.Range('I3').Value = 56
$ws.Range('I7').Value = 155

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I3').Value = 116
$ws.Range('I6').Value = 109
$ws.Range('I7').Value = 323

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('I4').Value = 2
$ws.Range('I6').Value = 16
$ws.Range('I7').Value = 61

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I2').Value = 89
$ws.Range('I5').Value = 26
$ws.Range('I8').Value = 521
$ws.Range('I13').Value = 12
$ws.Range('I14').Value = 41
$ws.Range('G20').Value = 620
$ws.Range('I20').Value = 221
$ws.Range('I23').Value = 74
$ws.Range('I29').Value = 537
$ws.Range('I30').Value = 26
$ws.Range('I33').Value = 380
$ws.Range('I37').Value = 265
$ws.Range('I42').Value = 287
$ws.Range('I43').Value = 72
$ws.Range('H48').Value = 333
$ws.Range('I48').Value = 90
$ws.Range('I49').Value = 55
$ws.Range('I51').Value = 74
$ws.Range('I52').Value = 166
$ws.Range('I53').Value = 84
$ws.Range('I54').Value = 188
$ws.Range('I55').Value = 92
$ws.Range('I67').Value = 323
$ws.Range('I68').Value = 29
$ws.Range('I73').Value = 75
$ws.Range('I76').Value = 128
$ws.Range('I78').Value = 110
$ws.Range('I79').Value = 216
$ws.Range('I83').Value = 160
$ws.Range('I84').Value = 61
$ws.Range('I85').Value = 385
$ws.Range('I87').Value = 13
$ws.Range('I88').Value = 72
$ws.Range('I89').Value = 87
$ws.Range('I91').Value = 96
$ws.Range('I92').Value = 25
$ws.Range('I93').Value = 50
$ws.Range('I94').Value = 72
$ws.Range('I98').Value = 53
$ws.Range('I99').Value = 155
$ws.Range('I100').Value = 12
$ws.Range('G101').Value = 24652
$ws.Range('H101').Value = 25971
$ws.Range('I101').Value = 8317

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I3').Value = 64
$ws.Range('I7').Value = 160

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I3').Value = 132
$ws.Range('I7').Value = 380

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('I6').Value = 31
$ws.Range('I7').Value = 55

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I6').Value = 93
$ws.Range('I7').Value = 188

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 169
$ws.Range('I3').Value = 182
$ws.Range('I6').Value = 153
$ws.Range('I7').Value = 537

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I3').Value = 18
$ws.Range('H4').Value = 41
$ws.Range('I6').Value = 52
$ws.Range('H7').Value = 333
$ws.Range('I7').Value = 90

$ws = $wb.Worksheets.Item('River North')
$ws.Range('I4').Value = 19
$ws.Range('I6').Value = 54
$ws.Range('I7').Value = 128

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I2').Value = 94
$ws.Range('I3').Value = 151
$ws.Range('I6').Value = 111
$ws.Range('I7').Value = 385

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I2').Value = 74
$ws.Range('I4').Value = 27
$ws.Range('I7').Value = 287

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('I3').Value = 3
$ws.Range('I6').Value = 12

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I3').Value = 30
$ws.Range('I7').Value = 110

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('I2').Value = 32
$ws.Range('I7').Value = 92

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('I3').Value = 27
$ws.Range('I7').Value = 74

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('I6').Value = 28
$ws.Range('I7').Value = 96

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I3').Value = 65
$ws.Range('I6').Value = 73
$ws.Range('I7').Value = 216

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I2').Value = 57
$ws.Range('I3').Value = 66
$ws.Range('G4').Value = 26
$ws.Range('I4').Value = 13
$ws.Range('I6').Value = 79
$ws.Range('G7').Value = 620
$ws.Range('I7').Value = 221

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('I3').Value = 15
$ws.Range('I7').Value = 50

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('I3').Value = 5
$ws.Range('I6').Value = 12

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I2').Value = 41
$ws.Range('I3').Value = 67
$ws.Range('I7').Value = 166

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I6').Value = 42
$ws.Range('I7').Value = 72

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('I6').Value = 34
$ws.Range('I7').Value = 53

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I3').Value = 25
$ws.Range('I6').Value = 20
$ws.Range('I7').Value = 75

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('I2').Value = 27
$ws.Range('I3').Value = 33
$ws.Range('I6').Value = 18
$ws.Range('I7').Value = 89

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('I2').Value = 9
$ws.Range('I7').Value = 25

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I6').Value = 27
$ws.Range('I7').Value = 72

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 163
$ws.Range('I3').Value = 139
$ws.Range('I4').Value = 32
$ws.Range('I6').Value = 171
$ws.Range('I7').Value = 521

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('I3').Value = 6
$ws.Range('I7').Value = 26

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I6').Value = 36
$ws.Range('I7').Value = 74

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('I4').Value = 5
$ws.Range('I7').Value = 29

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('I3').Value = 13
$ws.Range('I7').Value = 72

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('I6').Value = 36
$ws.Range('I7').Value = 84

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('I3').Value = 6
$ws.Range('I7').Value = 13
